# Applies the "added basal maturation rate" edit to
# rate_matrix_Maturation_sensor.xlsx (sheet ParallelMaturationAndPriming).
#
# Summary of the change:
#  - A new parameter row is inserted above the old "k_mature" row (old row 49),
#    pushing all subsequent parameter rows down by one.
#  - The new row holds a new named parameter "k_mature_basal" = 2.
#  - All the named ranges that pointed into the parameter block are
#    re-pointed to their (shifted) target cells, and the new name
#    "k_mature_basal" is created.
#  - C4 (basal maturation-rate cell of the base rate matrix) now refers to
#    k_mature_basal instead of k_mature.
#  - The Ca-dependent rate terms (k_prime, k_mature, and the various
#    "n*k_on" terms) are moved out of the base rate matrix (rows 2-22) and
#    into the "Ca Dependence" matrix (rows 24-45), which previously just
#    held binary 0/1 flags.
#  - The "Ca Dependence (binary)" header becomes "Ca Dependence" and the
#    header row is taller to fit the new explanatory text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ParallelMaturationAndPriming")

# ---------------------------------------------------------------------
# 1. Insert a new row above the old "k_mature" parameter row (row 49),
#    shifting the rest of the parameter block down by one row.
# ---------------------------------------------------------------------
$ws.Rows("49:49").Insert()

# The row that used to be row 49 (now row 50, "k_mature") carried a custom
# row height down from what is now row 51 ("k_unmature", which had a
# custom height of 30 for its old wrapped label). The final layout has no
# custom height on row 51, so reset it back to the sheet default.
$ws.Rows("51:51").AutoFit()

# ---------------------------------------------------------------------
# 2. Re-point the existing named ranges to their new (shifted) rows, and
#    create the new "k_mature_basal" name.
# ---------------------------------------------------------------------
$wb.Names.Item("k_mature").RefersTo     = "=ParallelMaturationAndPriming!`$B`$50"
$wb.Names.Item("k_unmature").RefersTo   = "=ParallelMaturationAndPriming!`$B`$51"
$wb.Names.Item("k_prime").RefersTo      = "=ParallelMaturationAndPriming!`$B`$52"
$wb.Names.Item("k_unprime").RefersTo    = "=ParallelMaturationAndPriming!`$B`$53"
$wb.Names.Item("k_on").RefersTo         = "=ParallelMaturationAndPriming!`$B`$54"
$wb.Names.Item("k_off").RefersTo        = "=ParallelMaturationAndPriming!`$B`$55"
$wb.Names.Item("b").RefersTo            = "=ParallelMaturationAndPriming!`$B`$56"
$wb.Names.Item("f").RefersTo            = "=ParallelMaturationAndPriming!`$B`$57"
$wb.Names.Item("M_plus").RefersTo       = "=ParallelMaturationAndPriming!`$B`$58"
$wb.Names.Item("P_plus").RefersTo       = "=ParallelMaturationAndPriming!`$B`$59"
# k_refill stays on row 48, unchanged.

$wb.Names.Add("k_mature_basal", "=ParallelMaturationAndPriming!`$B`$49")

# ---------------------------------------------------------------------
# 3. Fill in the new parameter row (basal maturation rate).
#    Column B in this block is normally formatted as Text ("@"), but the
#    new value cell uses a plain General/left-aligned style (like the
#    rest of the workbook's numeric cells) so it is stored as a genuine
#    number rather than text.
# ---------------------------------------------------------------------
$ws.Range("A49").Value = "k_mature_basal"
$ws.Range("B49").NumberFormat = "General"
$ws.Range("B49").HorizontalAlignment = -4131
$ws.Range("B49").WrapText = $false
$ws.Range("B49").Value = 2

# ---------------------------------------------------------------------
# 4. Update the "Ca Dependence" matrix header (row 24).
# ---------------------------------------------------------------------
$ws.Range("A24").Value = "Ca Dependence"
$ws.Rows("24:24").RowHeight = 75

# ---------------------------------------------------------------------
# 5. The basal rate matrix (rows 2-22): the basal maturation cell (C4)
#    now references the new basal rate; the priming / Ca-dependent rate
#    cells lose their formulas (the rates move to the Ca-dependence
#    matrix below), leaving plain 0s in their place.
# ---------------------------------------------------------------------
$ws.Range("C4").Formula = "=k_mature_basal"

$ws.Range("C5").Value  = 0
$ws.Range("D7").Value  = 0
$ws.Range("E7").Value  = 0
$ws.Range("D8").Value  = 0
$ws.Range("H9").Value  = 0
$ws.Range("I10").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("N15").Value = 0
$ws.Range("O16").Value = 0
$ws.Range("P17").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("R19").Value = 0
$ws.Range("S20").Value = 0
$ws.Range("T21").Value = 0
$ws.Range("U22").Value = 0

# ---------------------------------------------------------------------
# 6. The Ca-dependence matrix (rows 24-45): the cells that used to hold a
#    plain "1" flag now hold the actual Ca-dependent rate formulas that
#    were removed from the basal matrix above.
# ---------------------------------------------------------------------
$ws.Range("C27").Formula = "=k_mature"
$ws.Range("C28").Formula = "=k_prime"
$ws.Range("D30").Formula = "=k_prime"
$ws.Range("E30").Formula = "=k_mature"
$ws.Range("D31").Formula = "=5*k_on"
$ws.Range("H32").Formula = "=4*k_on"
$ws.Range("I33").Formula = "=3*k_on"
$ws.Range("J34").Formula = "=2*k_on"
$ws.Range("K35").Formula = "=k_on"
$ws.Range("E36").Formula = "=5*k_on"
$ws.Range("M37").Formula = "=4*k_on"
$ws.Range("N38").Formula = "=3*k_on"
$ws.Range("O39").Formula = "=2*k_on"
$ws.Range("P40").Formula = "=k_on"
$ws.Range("G41").Formula = "=5*k_on"
$ws.Range("R42").Formula = "=4*k_on"
$ws.Range("S43").Formula = "=3*k_on"
$ws.Range("T44").Formula = "=2*k_on"
$ws.Range("U45").Formula = "=k_on"

# ---------------------------------------------------------------------
# 7. Refresh the view (matches author's saved selection state).
# ---------------------------------------------------------------------
$ws.Range("B48").Select()

$wb.Application.CalculateFull()
